# dmmd examen - 4 de 10
# Enter exam grading data (EXAMEN sheet) for several students.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EXAMEN")
$ws.Activate()

# Row 3 - ALONSO NICOLAS / MATTEO AVNDAÑO
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 2
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 2
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0

# Row 5 - ANIBAL IGOR / QUIROZ LOPEZ
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 2
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0

# Row 6 - BASTIAN ANTONIO / CELEDON AGUIRRE
$ws.Range("D6").Value = 1

# Row 7 - HINOJOSA BARRIGA
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 2
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 2
$ws.Range("N7").Value = 3
$ws.Range("O7").Value = 2
$ws.Range("P7").Value = 2
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = 2
$ws.Range("U7").Value = 0

# Row 8 - ESTAY COLISTRO
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 2
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 2
$ws.Range("N8").Value = 3
$ws.Range("O8").Value = 2
$ws.Range("P8").Value = 2
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 2
$ws.Range("U8").Value = 1

# Row 9 - BARROS AVILA
$ws.Range("D9").Value = 3

# Row 10 - LY CASTRO
$ws.Range("D10").Value = 5
$ws.Range("G10").HorizontalAlignment = -4131

# Row 14 - CHAPARRO PEÑA
$ws.Range("D14").Value = 2

# Row 18 - VALDES ALAMOS
$ws.Range("D18").Value = 4

# Apply the "0.0" number format used throughout the grade columns to the
# full computed-grade ranges (I, M, Q, V, W for rows 3:19), matching the
# formatting already used elsewhere in the workbook for these columns.
$ws.Range("I3:I19").NumberFormat = "0.0"
$ws.Range("M3:M19").NumberFormat = "0.0"
$ws.Range("Q3:Q19").NumberFormat = "0.0"
$ws.Range("V3:V19").NumberFormat = "0.0"
$ws.Range("W3:W19").NumberFormat = "0.0"

# Restore the active-cell selection as left by the author.
$ws.Range("D19").Select()

$excel.CalculateFull()
